# YumaPointImport.xlsx - Point work: fix the calculated "Longitude" column.
#
# The table's Longitude column (AB) was mistakenly built from the Latitude
# source column (J = PRIM_LAT_DEC) instead of the Longitude source column
# (K = PRIM_LONG_DEC). Correct the formula so Longitude pulls from K2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AB2").Formula = "=K2"

# Reflect the author's final cursor position/selection on the sheet.
$ws.Range("AB3").Select() | Out-Null
